$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.835.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.242.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '269.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0964'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +19.45%  '
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.574.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.230.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.805'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.821.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("E19").Value = '  +2.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.80%  '
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.92%  '
$ws.Range("E27").Value = '  +12.49%  '
$ws.Range("E28").Value = '  +5.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '41.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0919'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.16%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.124'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.114'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0350'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +25.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.228'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '63.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.22%  '
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("E48").Value = '  +4.51%  '
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.442'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.459.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.20%  '
